$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "Source" block at the bottom of the sheet used to read:
#   Source:
#   NUMBER OF BUSINESS ESTABLISHMENTS BY EMPLOYMENT SIZE AND INDUSTRY (TTSNA) - 2007
#   http://cso.planning.gov.tt/.../Number of Business Establishments....pdf   (hyperlinked)
#   <blank>
#   T&TCSO
#   Trinidad & Tobago Central Statistics Office (T&TCSO), Business Establishments ...
#   TTCIC
#   Trinidad & Tobago Chamber of Industry and Commerce (TTCIC), "Small And Medium..."
#
# It becomes:
#   Source:
#   <blank>
#   NUMBER OF BUSINESS ESTABLISHMENTS BY EMPLOYMENT SIZE AND INDUSTRY (TTSNA) - 2007
#   <blank>
#   http://cso.planning.gov.tt/.../Number of Business Establishments....pdf   (plain text, no hyperlink)
#   T&TCSO
#   T&TCSO
#   TTCIC
#   "Enterprise Development Policy and Strategic Plan 2001 - 2005," Ministry of
#   Enterprise Development and Foreign Affairs, July 2001, p. 35, ...
#
# i.e. the hyperlink + its citation text are dropped (the citation text row
# is reused, now just repeating the "T&TCSO" label) and the old TTCIC citation
# text is replaced with a brand new citation.
# ---------------------------------------------------------------------------

# Drop the hyperlink on the URL cell (A76) entirely.
$ws.Hyperlinks.Delete()

# Remove the old hyperlinked URL row outright (its text moves a couple of
# rows down, re-typed below) - this also guarantees none of the old
# hyperlink-flavoured styling survives on the cell.
$ws.Rows("76:76").Delete()

# Re-open a blank line above the "NUMBER OF BUSINESS..." title line.
$ws.Rows("75:75").Insert()

# Re-open a blank line below the title (for the re-typed URL line).
$ws.Rows("78:78").Insert()
$ws.Range("A78").Value = "http://cso.planning.gov.tt/sites/default/files/content/documents/Number%20of%20Business%20Establishments%20in%20Trinidad%20_%20Tobago%20by%20Employment%20Size%20and%20Industry%202007.pdf"

# The T&TCSO citation paragraph is gone; the row is re-used to just repeat
# the "T&TCSO" label.
$ws.Range("A82").Value = "T&TCSO"

# The TTCIC citation paragraph is replaced with a new source altogether.
$ws.Range("A84").Value = '"Enterprise Development Policy and Strategic Plan 2001 – 2005,"Ministry of Enterprise Development and Foreign Affairs, July 2001, p. 35, http://www.sice.oas.org/ctyindex/TTO/INDPolicy_e.pdf'
